# Applies the nBackCond3.xlsx update:
#  - rename header labels (drop "T3" suffix)
#  - update selection to B31
#  - update a handful of trial data values in columns A-C

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: strip "T3" suffix from column labels
# A1: squareT3 -> square, B1: loc1T3 -> loc1, C1: loc2T3 -> loc2, D1: corrAnsT3 -> corrAns
$ws.Range("A1").Value = "square"
$ws.Range("B1").Value = "loc1"
$ws.Range("C1").Value = "loc2"
$ws.Range("D1").Value = "corrAns"

$ws.Range("A2").Value = 4
$ws.Range("B2").Value = -0.2

$ws.Range("B4").Value = 0.2

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = -0.2
$ws.Range("C5").Value = 0

$ws.Range("A8").Value = 1
$ws.Range("C8").Value = 0.2

$ws.Range("A11").Value = 1
$ws.Range("C11").Value = 0.2

$ws.Range("A16").Value = 7
$ws.Range("B16").Value = -0.2

$ws.Range("A20").Value = 7
$ws.Range("B20").Value = -0.2
$ws.Range("C20").Value = -0.2

$ws.Range("A23").Value = 8
$ws.Range("B23").Value = 0
$ws.Range("C23").Value = -0.2

$ws.Range("A26").Value = 8
$ws.Range("B26").Value = 0
$ws.Range("C26").Value = -0.2

$ws.Range("B28").Value = 0
$ws.Range("C28").Value = 0.2

$ws.Range("B31").Value = 0
$ws.Range("C31").Value = 0.2

# Update selection to B31 (matches the saved selection in the diff)
$ws.Range("B31").Select()
